$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "B2" = 8.292000000000007
    "E2" = 16.66670000000001
    "A3" = -21.6394
    "D3" = -7.275599999999995
    "E6" = 16.50910000000001
    "D12" = -7.336800000000002
    "A14" = -21.66629999999999
    "A16" = -21.74339999999999
    "B18" = 6.805899999999999
    "E19" = 16.2626
    "A21" = -20.25079999999997
    "A23" = -20.81319999999997
    "B24" = 5.516800000000002
    "D24" = -7.757300000000002
    "E24" = 16.7023
    "A25" = -21.85199999999999
    "B25" = 6.443600000000004
    "D25" = -8.984499999999992
    "A26" = -21.21809999999997
    "B27" = 6.092800000000003
    "E27" = 16.65709999999999
    "A29" = -20.70969999999997
    "B30" = 6.571200000000001
    "E30" = 15.43529999999999
    "B31" = 6.207800000000004
    "E31" = 15.87629999999999
    "E33" = 16.86620000000002
    "B39" = 9.286000000000007
    "A40" = -19.4483
    "D41" = -8.569499999999993
    "B42" = 10.05659999999999
    "E42" = 16.4403
    "B48" = 5.7057
    "D50" = -8.066500000000001
    "B51" = 5.730799999999999
    "B52" = 5.397499999999999
    "A53" = -21.9583
    "D53" = -6.081300000000002
    "B55" = 6.566299999999996
    "E55" = 16.511
    "B56" = 5.606299999999997
    "D56" = -7.865099999999995
    "A57" = -22.34340000000001
    "B57" = 4.405599999999994
    "D57" = -8.555600000000005
    "D58" = -8.504700000000007
    "E58" = 16.3822
    "A59" = -22.4848
    "B60" = 5.659699999999999
    "D61" = -7.6323
    "D63" = -7.598200000000004
    "D64" = -7.501299999999998
    "A65" = -21.88249999999998
    "E65" = 17.19150000000001
    "A69" = -21.5896
    "D70" = -8.029800000000005
    "E70" = 16.5714
    "D72" = -7.615299999999998
    "B73" = 8.721399999999997
    "B74" = 9.427499999999988
    "E74" = 16.36579999999999
    "E75" = 16.66730000000001
    "A79" = -20.447
    "A83" = -21.9218
    "E83" = 16.81250000000001
    "E84" = 16.34219999999999
    "D86" = -8.572500000000003
    "E86" = 16.32990000000001
    "B89" = 4.778599999999992
    "D89" = -5.728000000000003
    "B90" = 6.203900000000002
    "A91" = -21.35310000000001
    "B92" = 4.926699999999993
    "A93" = -21.02979999999998
    "E96" = 16.46769999999999
    "E97" = 16.89150000000001
    "D98" = -8.691899999999999
    "A100" = -22.03649999999999
    "D100" = -8.329999999999998
    "D102" = -7.987199999999996
}

foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}
